$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the "_GoBack" bookmark that sits right after the sentence
# "L'utilisateur peut checker la prise de médicament".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: turn "MyTherapie" into two runs "MyTherap" + "y" (i.e. correct the
# spelling to "MyTherapy") and re-insert the "_GoBack" bookmark right after
# the newly typed "y", i.e. at the very place the user last edited.
# ---------------------------------------------------------------------------
$oldWord = "MyTherapie"
$newWord = "MyTherapy"
$keepLen = 8   # length of the common prefix "MyTherap"

$rng = $d.Content
$found = $rng.Find.Execute($oldWord, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $startPos = $rng.Start
    $endPos = $rng.End

    # Replace the differing tail ("ie") with the new tail ("y"). This keeps
    # everything inside the single original run for now ("MyTherapy").
    $oldTail = $oldWord.Substring($keepLen)
    $newTail = $newWord.Substring($keepLen)
    $tailRange = $d.Range($endPos - $oldTail.Length, $endPos)
    $tailRange.Text = $newTail

    $newEndPos = $startPos + $newWord.Length

    # Force a clean run split between "MyTherap" and "y" by briefly adding
    # and removing a bookmark at that boundary (adding text directly at a
    # run boundary here would instead merge the surrounding runs).
    $splitPos = $startPos + $keepLen
    $splitPoint = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("TempSplitMarker", $splitPoint) | Out-Null
    $d.Bookmarks.Item("TempSplitMarker").Delete()

    # Re-create the "_GoBack" bookmark immediately after the "y", i.e. right
    # before the closing <w:proofErr .../> of the word. Inserting a bookmark
    # exactly at that spot directly is unreliable, so a temporary character
    # is typed after the target spot, the bookmark is anchored just before
    # it, and the temporary character is removed again.
    $afterWord = $newEndPos
    $tempPoint = $d.Range($afterWord, $afterWord)
    $tempPoint.InsertAfter("Z")
    $bmPoint = $d.Range($afterWord, $afterWord)
    $d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null
    $tempCharRange = $d.Range($afterWord, $afterWord + 1)
    $tempCharRange.Text = ""
}
